$wb = $excel.ActiveWorkbook

# Rename sheets (figures/tables renumbering: S1-S7 -> A-G, S8/S9 drop "Table S#" prefix)
$wb.Worksheets.Item("Table S1 - Sample Size").Name = "Table A - Sample Size"
$wb.Worksheets.Item("Table S2 - PERMANOVA AIC").Name = "Table B - PERMANOVA AIC"
$wb.Worksheets.Item("Table S3 - Plasticity AIC").Name = "Table C - Plasticity AIC"
$wb.Worksheets.Item("Table S4 - PERMANOVA").Name = "Table D - PERMANOVA"
$wb.Worksheets.Item("Table S5 - Plasticity GLM").Name = "Table E - Plasticity GLM"
$wb.Worksheets.Item("Table S6 - Species PERMANOVA").Name = "Table F - Species PERMANOVA"
$wb.Worksheets.Item("Table S7 - HostVsymb PERMANOVA").Name = "Table G - HostVsymb PERMANOVA"
$wb.Worksheets.Item("Table S8 - HostVsymb Plast AIC").Name = "HostVsymb Plast AIC"
$wb.Worksheets.Item("Table S9 - HostVsymb Plast GLM").Name = "HostVsymb Plast GLM"

# Update p-values on "Table F - Species PERMANOVA" (formerly Table S6)
$wsF = $wb.Worksheets.Item("Table F - Species PERMANOVA")
$wsF.Range("F3").Value = 0.09927
$wsF.Range("F4").Value = 0.002
$wsF.Range("F7").Value = 0.01932
$wsF.Range("F8").Value = 0.002

# Update p-values on "Table G - HostVsymb PERMANOVA" (formerly Table S7)
$wsG = $wb.Worksheets.Item("Table G - HostVsymb PERMANOVA")
$wsG.Range("F2").Value = 0.77815
$wsG.Range("F4").Value = 0.59227
$wsG.Range("F7").Value = 0.26316
$wsG.Range("F9").Value = 0.17055
$wsG.Range("F12").Value = 0.01732
$wsG.Range("F14").Value = 0.1972
